# Actualiza datos y proyeccion inicial GAMLSS
# Adds the 2022 population data point (a first GAMLSS-based projection) as a
# new last row, moving the "final row" emphasis formatting (bold, bordered,
# right aligned number format) from the previous last row (2021, row 23) to
# the newly appended last row (2022, row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Give the new last row (B24) the same "final row" formatting that
#        currently lives on B23 (bold / bordered / right aligned), BEFORE we
#        touch B23's own formatting.
$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. Re-format the old last row (B23) back to the regular data-row style
#        (same look as B2:B22), since it is no longer the final row.
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)  # xlPasteFormats

# --- 3. Write the new row's values (2022 / 5213362).
$ws.Range("A24").Value = 2022
$ws.Range("B24").Value = 5213362

# --- 4. Reset the active selection to A1, like the saved workbook shows.
$ws.Range("A1").Select()
